$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.915.47"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.362.34"
$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("D5").Value = "'302.50"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").Value = "'95.73"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "'34.06"
$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("E11").Value = "  +3.45%  "

$ws.Range("D12").Value = "'0.0784"

$ws.Range("D13").Value = "'18.28"
$ws.Range("E13").Value = "  -2.91%  "

$ws.Range("D14").Value = "'6.72"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("D15").Value = "2.729.31"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").Value = "2.349.38"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "42.882.08"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").Value = "'11.82"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").Value = "0.0₃0883"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").Value = "'67.89"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "'235.02"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  -4.98%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "'24.49"
$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("D29").Value = "'9.30"
$ws.Range("E29").Value = "  +2.27%  "

$ws.Range("D30").Value = "'31.85"

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").Value = "'5.01"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").Value = "'17.34"
$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("D34").Value = "'0.0710"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'127.89"
$ws.Range("E35").Value = "  -22.98%  "

$ws.Range("E36").Value = "  +3.09%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.104"
$ws.Range("E37").Value = "  +3.52%  "

$ws.Range("E38").Value = "  -2.81%  "

$ws.Range("E39").Value = "  -2.32%  "

$ws.Range("E40").Value = "  +3.14%  "

$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("D42").Value = "'21.05"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").Value = "1.928.12"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.16"
$ws.Range("E47").Value = "  -8.65%  "

$ws.Range("D48").Value = "2.588.92"
$ws.Range("E48").Value = "  +1.94%  "

$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D51").Value = "'51.53"
$ws.Range("E51").Value = "  -3.31%  "
